$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 84
$ws.Range("D2").Value = 28
$ws.Range("E2").Value = 56
$ws.Range("F2").Value = 54.90196078431373

$ws.Range("C3").Value = 194
$ws.Range("D3").Value = 47
$ws.Range("E3").Value = 147
$ws.Range("F3").Value = 55.95238095238096

$ws.Range("C4").Value = 186
$ws.Range("D4").Value = 81
$ws.Range("E4").Value = 105
$ws.Range("F4").Value = 41.75257731958763

$ws.Range("C5").Value = 210
$ws.Range("D5").Value = 94
$ws.Range("E5").Value = 116
$ws.Range("F5").Value = 50.53763440860215

$ws.Range("C6").Value = 216
$ws.Range("D6").Value = 128
$ws.Range("E6").Value = 88
$ws.Range("F6").Value = 60.95238095238096

$ws.Range("C7").Value = 235
$ws.Range("D7").Value = 129
$ws.Range("E7").Value = 106
$ws.Range("F7").Value = 59.72222222222222

$ws.Range("C8").Value = 191
$ws.Range("D8").Value = 130
$ws.Range("E8").Value = 61
$ws.Range("F8").Value = 55.31914893617022

$ws.Range("C9").Value = 174
$ws.Range("D9").Value = 126
$ws.Range("E9").Value = 48
$ws.Range("F9").Value = 65.96858638743456

$ws.Range("C10").Value = 164
$ws.Range("D10").Value = 124
$ws.Range("E10").Value = 40
$ws.Range("F10").Value = 71.26436781609196

$ws.Range("C11").Value = 194
$ws.Range("D11").Value = 142
$ws.Range("E11").Value = 52
$ws.Range("F11").Value = 86.58536585365853

$ws.Range("C12").Value = 200
$ws.Range("D12").Value = 152
$ws.Range("E12").Value = 48
$ws.Range("F12").Value = 78.35051546391753

$ws.Range("C13").Value = 213
$ws.Range("D13").Value = 157
$ws.Range("E13").Value = 56
$ws.Range("F13").Value = 78.5

$ws.Range("C14").Value = 222
$ws.Range("D14").Value = 181
$ws.Range("E14").Value = 41
$ws.Range("F14").Value = 84.97652582159625

$ws.Range("C15").Value = 247
$ws.Range("D15").Value = 181
$ws.Range("E15").Value = 66
$ws.Range("F15").Value = 81.53153153153153

$ws.Range("C16").Value = 310
$ws.Range("D16").Value = 204
$ws.Range("E16").Value = 106
$ws.Range("F16").Value = 82.5910931174089

$ws.Range("C17").Value = 314
$ws.Range("D17").Value = 248
$ws.Range("E17").Value = 66
$ws.Range("F17").Value = 80

$ws.Range("C18").Value = 310
$ws.Range("D18").Value = 250
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 79.61783439490446

$ws.Range("C19").Value = 304
$ws.Range("D19").Value = 248
$ws.Range("E19").Value = 56
$ws.Range("F19").Value = 80

$ws.Range("C20").Value = 237
$ws.Range("D20").Value = 206
$ws.Range("E20").Value = 31
$ws.Range("F20").Value = 67.76315789473685
